$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 36000.1
$ws.Cells.Item(64, 9).Value = 252050
$ws.Cells.Item(64, 10).Value = 2761.6538
$ws.Cells.Item(64, 11).Value = 252050
$ws.Cells.Item(64, 12).Value = 2761.6538
$ws.Cells.Item(64, 13).Value = -251802
$ws.Cells.Item(64, 14).Value = -3257.6538
$ws.Cells.Item(67, 8).Value = 36000.1
$ws.Cells.Item(67, 9).Value = 252050
$ws.Cells.Item(67, 10).Value = 2761.6538
$ws.Cells.Item(67, 11).Value = 252050
$ws.Cells.Item(67, 12).Value = 2761.6538
$ws.Cells.Item(67, 13).Value = -251192
$ws.Cells.Item(67, 14).Value = -4477.6538
$ws.Cells.Item(95, 8).Value = 33205.332
$ws.Cells.Item(95, 10).Value = 33205.332
$ws.Cells.Item(95, 12).Value = 33205.332
$ws.Cells.Item(95, 14).Value = -38697.332
$ws.Cells.Item(117, 8).Value = 48333.5
$ws.Cells.Item(117, 10).Value = 48333.5
$ws.Cells.Item(117, 12).Value = 48333.5
$ws.Cells.Item(117, 14).Value = -57511.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 51400.57
$ws.Cells.Item(80, 10).Value = 51400.57
$ws.Cells.Item(80, 12).Value = 51400.57
$ws.Cells.Item(80, 14).Value = -53396.57
$ws.Cells.Item(83, 8).Value = 51400.57
$ws.Cells.Item(83, 10).Value = 51400.57
$ws.Cells.Item(83, 12).Value = 154201.71
$ws.Cells.Item(83, 14).Value = -164185.71
$ws.Cells.Item(104, 8).Value = 41092.5
$ws.Cells.Item(104, 10).Value = 41092.5
$ws.Cells.Item(104, 12).Value = 41092.5
$ws.Cells.Item(104, 14).Value = -48080.5
$ws.Cells.Item(105, 8).Value = 47945.8
$ws.Cells.Item(105, 10).Value = 47945.8
$ws.Cells.Item(105, 12).Value = 47945.8
$ws.Cells.Item(105, 14).Value = -54933.8
$ws.Cells.Item(106, 8).Value = 46870.8
$ws.Cells.Item(106, 10).Value = 46870.8
$ws.Cells.Item(106, 12).Value = 46870.8
$ws.Cells.Item(106, 14).Value = -49394.8
$ws.Cells.Item(107, 8).Value = 36257
$ws.Cells.Item(107, 10).Value = 36257
$ws.Cells.Item(107, 12).Value = 36257
$ws.Cells.Item(107, 14).Value = -43937
$ws.Cells.Item(109, 8).Value = 43151
$ws.Cells.Item(109, 10).Value = 43151
$ws.Cells.Item(109, 12).Value = 43151
$ws.Cells.Item(109, 14).Value = -45925
$ws.Cells.Item(117, 8).Value = 47918.5
$ws.Cells.Item(117, 10).Value = 47918.5
$ws.Cells.Item(117, 12).Value = 47918.5
$ws.Cells.Item(117, 14).Value = -57096.5
$ws.Cells.Item(118, 8).Value = 49366.668
$ws.Cells.Item(118, 10).Value = 49366.668
$ws.Cells.Item(118, 12).Value = 49366.668
$ws.Cells.Item(118, 14).Value = -52680.668
$ws.Cells.Item(119, 8).Value = 52588
$ws.Cells.Item(119, 10).Value = 52588
$ws.Cells.Item(119, 12).Value = 52588
$ws.Cells.Item(119, 14).Value = -62264
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(117, 8).Value = 49911.332
$ws.Cells.Item(117, 10).Value = 49911.332
$ws.Cells.Item(117, 12).Value = 49911.332
$ws.Cells.Item(117, 14).Value = -59089.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(92, 8).Value = 36381.082
$ws.Cells.Item(92, 10).Value = 36381.082
$ws.Cells.Item(92, 12).Value = 36381.082
$ws.Cells.Item(92, 14).Value = -41373.082
$ws.Cells.Item(104, 8).Value = 31077.625
$ws.Cells.Item(104, 10).Value = 31077.625
$ws.Cells.Item(104, 12).Value = 31077.625
$ws.Cells.Item(104, 14).Value = -36319.625
$ws.Cells.Item(109, 8).Value = 28329.3
$ws.Cells.Item(109, 10).Value = 28329.3
$ws.Cells.Item(109, 12).Value = 28329.3
$ws.Cells.Item(109, 14).Value = -30409.3
$ws.Cells.Item(111, 8).Value = 47264
$ws.Cells.Item(111, 10).Value = 47264
$ws.Cells.Item(111, 12).Value = 47264
$ws.Cells.Item(111, 14).Value = -55444
$ws.Cells.Item(115, 8).Value = 30788.4
$ws.Cells.Item(115, 10).Value = 30788.4
$ws.Cells.Item(115, 12).Value = 30788.4
$ws.Cells.Item(115, 14).Value = -33138.4
$ws.Cells.Item(116, 8).Value = 47785
$ws.Cells.Item(116, 10).Value = 47785
$ws.Cells.Item(116, 12).Value = 47785
$ws.Cells.Item(116, 14).Value = -56963
$ws.Cells.Item(120, 8).Value = 31822.416
$ws.Cells.Item(120, 10).Value = 31822.416
$ws.Cells.Item(120, 12).Value = 31822.416
$ws.Cells.Item(120, 14).Value = -39080.416
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2490.291
$ws.Cells.Item(113, 9).Value = 3989.2334
$ws.Cells.Item(113, 10).Value = 691.5599999999999
$ws.Cells.Item(113, 11).Value = 11967.7002
$ws.Cells.Item(113, 12).Value = 2074.68
$ws.Cells.Item(113, 13).Value = -9797.700199999999
$ws.Cells.Item(113, 14).Value = -6414.68
$ws.Cells.Item(133, 8).Value = 8099.778
$ws.Cells.Item(133, 9).Value = 10233.333
$ws.Cells.Item(133, 10).Value = 7033
$ws.Cells.Item(133, 11).Value = 30699.999
$ws.Cells.Item(133, 12).Value = 21099
$ws.Cells.Item(133, 13).Value = -25639.999
$ws.Cells.Item(133, 14).Value = -31219
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 44919
$ws.Cells.Item(104, 10).Value = 44919
$ws.Cells.Item(104, 12).Value = 44919
$ws.Cells.Item(104, 14).Value = -51907
$ws.Cells.Item(105, 8).Value = 42940.5
$ws.Cells.Item(105, 10).Value = 42940.5
$ws.Cells.Item(105, 12).Value = 42940.5
$ws.Cells.Item(105, 14).Value = -49928.5
$ws.Cells.Item(114, 8).Value = 40643.285
$ws.Cells.Item(114, 10).Value = 40643.285
$ws.Cells.Item(114, 12).Value = 40643.285
$ws.Cells.Item(114, 14).Value = -49321.285
$ws.Cells.Item(116, 8).Value = 38938
$ws.Cells.Item(116, 10).Value = 38938
$ws.Cells.Item(116, 12).Value = 38938
$ws.Cells.Item(116, 14).Value = -48116
$ws.Cells.Item(118, 8).Value = 38198.332
$ws.Cells.Item(118, 10).Value = 38198.332
$ws.Cells.Item(118, 12).Value = 38198.332
$ws.Cells.Item(118, 14).Value = -41512.332
$ws.Cells.Item(122, 8).Value = 1173.75
$ws.Cells.Item(122, 9).Value = 1231.6666
$ws.Cells.Item(122, 11).Value = 3694.9998
$ws.Cells.Item(122, 13).Value = -1244.9998
$ws.Cells.Item(126, 8).Value = 3622.4
$ws.Cells.Item(126, 9).Value = 5756
$ws.Cells.Item(126, 11).Value = 17268
$ws.Cells.Item(126, 13).Value = -14798
$ws.Cells.Item(130, 8).Value = 44425.332
$ws.Cells.Item(130, 10).Value = 44425.332
$ws.Cells.Item(130, 12).Value = 44425.332
$ws.Cells.Item(130, 14).Value = -54465.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(97, 8).Value = 33500
$ws.Cells.Item(97, 10).Value = 33500
$ws.Cells.Item(97, 12).Value = 33500
$ws.Cells.Item(97, 14).Value = -35482
$ws.Cells.Item(100, 8).Value = 2145.0908
$ws.Cells.Item(100, 9).Value = 2059.6
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 2059.6
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = -1518.6
$ws.Cells.Item(100, 14).Value = -4082
$ws.Cells.Item(103, 8).Value = 40824.8
$ws.Cells.Item(103, 10).Value = 40824.8
$ws.Cells.Item(103, 12).Value = 40824.8
$ws.Cells.Item(103, 14).Value = -43168.8
$ws.Cells.Item(110, 8).Value = 45544
$ws.Cells.Item(110, 10).Value = 45544
$ws.Cells.Item(110, 12).Value = 45544
$ws.Cells.Item(110, 14).Value = -53724
$ws.Cells.Item(129, 8).Value = 43425
$ws.Cells.Item(129, 10).Value = 43425
$ws.Cells.Item(129, 12).Value = 43425
$ws.Cells.Item(129, 14).Value = -53425
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 42336
$ws.Cells.Item(95, 10).Value = 42336
$ws.Cells.Item(95, 12).Value = 42336
$ws.Cells.Item(95, 14).Value = -47828
$ws.Cells.Item(105, 8).Value = 50056.25
$ws.Cells.Item(105, 10).Value = 50056.25
$ws.Cells.Item(105, 12).Value = 50056.25
$ws.Cells.Item(105, 14).Value = -57044.25
$ws.Cells.Item(129, 8).Value = 39429
$ws.Cells.Item(129, 10).Value = 39429
$ws.Cells.Item(129, 12).Value = 39429
$ws.Cells.Item(129, 14).Value = -49429
